# Updated cryptos list on Fri Nov  3 03:22:26 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.882.13"
$ws.Range("E2").Value = "  -2.86%  "
$ws.Range("D3").Value = "1.808.72"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.318"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0679"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0994"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "2.069.11"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("D13").Value = "1.812.67"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.660"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.28%  "
$ws.Range("D17").Value = "34.822.35"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.73%  "
$ws.Range("E29").Value = "  +6.35%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0546"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.81%  "
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.680"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "91.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.19%  "
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").Value = "1.306.08"
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.957"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.44%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -13.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.31%  "
$ws.Range("E46").Value = "  -3.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0509"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "1.990.35"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("E49").Value = "  +7.39%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.59%  "